$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1 (mirror style/format of existing header cell H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data for columns I and J, rows 2-9
$values = @{
    2 = @(1, 7)
    3 = @(1, 6)
    4 = @(1, 4)
    5 = @(1, 4)
    6 = @(8, 8)
    7 = @(5, 6)
    8 = @(4, 4)
    9 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
